$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine rows 2-6 into a single Python-tuple-style string in A2
$ws.Range("A2").Value = "('Rukh Egg', ['{3}{R}', 'Creature " + [char]0x2014 + " Bird Egg', 'When Rukh Egg dies, create a 4/4 red Bird creature token with flying at the beginning of the next end step.', '0/3'])"

# Remove the now-redundant rows 3-6 so the sheet's used range shrinks to A1:A2
$ws.Range("A3:A6").EntireRow.Delete()
